$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2142.8572
$ws.Range("J40").Value = 2785.7144
$ws.Range("L40").Value = 2785.7144
$ws.Range("N40").Value = -3135.7144

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 898.3333
$ws.Range("I43").Value = 900
$ws.Range("J43").Value = 895
$ws.Range("K43").Value = 900
$ws.Range("L43").Value = 895
$ws.Range("M43").Value = -831
$ws.Range("N43").Value = -1033

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4124.0586
$ws.Range("I76").Value = 4019
$ws.Range("K76").Value = 4019
$ws.Range("M76").Value = -3704

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4124.0586
$ws.Range("I79").Value = 4019
$ws.Range("K79").Value = 4019
$ws.Range("M79").Value = -2927

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2443.75
$ws.Range("I86").Value = 1950
$ws.Range("J86").Value = 2937.5
$ws.Range("K86").Value = 1950
$ws.Range("L86").Value = 2937.5
$ws.Range("M86").Value = -827
$ws.Range("N86").Value = -5183.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2443.75
$ws.Range("I89").Value = 1950
$ws.Range("J89").Value = 2937.5
$ws.Range("K89").Value = 9750
$ws.Range("L89").Value = 14687.5
$ws.Range("M89").Value = -4134
$ws.Range("N89").Value = -25919.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 8211
$ws.Range("I98").Value = 6260.2354
$ws.Range("J98").Value = 16501.75
$ws.Range("K98").Value = 6260.2354
$ws.Range("L98").Value = 16501.75
$ws.Range("M98").Value = -4762.2354
$ws.Range("N98").Value = -19497.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2198.2
$ws.Range("I112").Value = 900.25
$ws.Range("J112").Value = 2324.8293
$ws.Range("K112").Value = 2700.75
$ws.Range("L112").Value = 6974.4879
$ws.Range("M112").Value = -1592.75
$ws.Range("N112").Value = -9190.4879

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 8211
$ws.Range("I122").Value = 6260.2354
$ws.Range("J122").Value = 16501.75
$ws.Range("K122").Value = 18780.7062
$ws.Range("L122").Value = 49505.25
$ws.Range("M122").Value = -16330.7062
$ws.Range("N122").Value = -54405.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1469.4038
$ws.Range("I137").Value = 1133.5
$ws.Range("J137").Value = 2880.2
$ws.Range("K137").Value = 3400.5
$ws.Range("L137").Value = 8640.599999999999
$ws.Range("M137").Value = -850.5
$ws.Range("N137").Value = -13740.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2636068
$ws.Range("I138").Value = 7695068
$ws.Range("J138").Value = 5387.98
$ws.Range("K138").Value = 23085204
$ws.Range("L138").Value = 16163.94
$ws.Range("M138").Value = -23080064
$ws.Range("N138").Value = -26443.94

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1991.8889
$ws.Range("I61").Value = 1961.375
$ws.Range("J61").Value = 2236
$ws.Range("K61").Value = 1961.375
$ws.Range("L61").Value = 2236
$ws.Range("M61").Value = -1749.375
$ws.Range("N61").Value = -2660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3001.4
$ws.Range("J88").Value = 3602.3333
$ws.Range("L88").Value = 3602.3333
$ws.Range("N88").Value = -4414.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3001.4
$ws.Range("J91").Value = 3602.3333
$ws.Range("L91").Value = 3602.3333
$ws.Range("N91").Value = -6410.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1991.8889
$ws.Range("I136").Value = 1961.375
$ws.Range("J136").Value = 2236
$ws.Range("K136").Value = 5884.125
$ws.Range("L136").Value = 6708
$ws.Range("M136").Value = -3334.125
$ws.Range("N136").Value = -11808

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 33701.062
$ws.Range("I20").Value = 49293.19
$ws.Range("J20").Value = 3934.2727
$ws.Range("K20").Value = 49293.19
$ws.Range("L20").Value = 3934.2727
$ws.Range("M20").Value = -49046.19
$ws.Range("N20").Value = -4428.2727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1372.8572
$ws.Range("I99").Value = 922
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 922
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = 576
$ws.Range("N99").Value = -5496

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1254467.5
$ws.Range("I134").Value = 2226053.2
$ws.Range("J134").Value = 5285.7144
$ws.Range("K134").Value = 6678159.600000001
$ws.Range("L134").Value = 15857.1432
$ws.Range("M134").Value = -6675624.600000001
$ws.Range("N134").Value = -20927.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2590.0876
$ws.Range("I31").Value = 1389.6285
$ws.Range("J31").Value = 4499.909
$ws.Range("K31").Value = 1389.6285
$ws.Range("L31").Value = 4499.909
$ws.Range("M31").Value = -1094.6285
$ws.Range("N31").Value = -5089.909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2590.0876
$ws.Range("I34").Value = 1389.6285
$ws.Range("J34").Value = 4499.909
$ws.Range("K34").Value = 1389.6285
$ws.Range("L34").Value = 4499.909
$ws.Range("M34").Value = -1187.6285
$ws.Range("N34").Value = -4903.909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2276
$ws.Range("I105").Value = 1793.3334
$ws.Range("K105").Value = 1793.3334
$ws.Range("M105").Value = -46.33339999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1948.1428
$ws.Range("I132").Value = 1543.2333
$ws.Range("J132").Value = 4377.6
$ws.Range("K132").Value = 4629.699900000001
$ws.Range("L132").Value = 13132.8
$ws.Range("M132").Value = -2099.699900000001
$ws.Range("N132").Value = -18192.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1073584.4
$ws.Range("I12").Value = 28
$ws.Range("K12").Value = 84
$ws.Range("M12").Value = 89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2050.5
$ws.Range("I132").Value = 1300.6666
$ws.Range("K132").Value = 11705.9994
$ws.Range("M132").Value = -9175.999400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 100000
$ws.Range("J42").Value = 100000
$ws.Range("L42").Value = 100000
$ws.Range("N42").Value = -100970

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 36000
$ws.Range("J82").Value = 36000
$ws.Range("L82").Value = 36000
$ws.Range("N82").Value = -36766

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H85").Value = 36000
$ws.Range("J85").Value = 36000
$ws.Range("L85").Value = 36000
$ws.Range("N85").Value = -38652

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 18164.645
$ws.Range("I97").Value = 23669.783
$ws.Range("K97").Value = 23669.783
$ws.Range("M97").Value = -23173.783

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 40000
$ws.Range("J106").Value = 40000
$ws.Range("L106").Value = 40000
$ws.Range("N106").Value = -42524

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 34500
$ws.Range("J114").Value = 34500
$ws.Range("L114").Value = 34500
$ws.Range("N114").Value = -43178

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H115").Value = 100000
$ws.Range("J115").Value = 100000
$ws.Range("L115").Value = 100000
$ws.Range("N115").Value = -102350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3314.2856
$ws.Range("I126").Value = 2200
$ws.Range("J126").Value = 3618.182
$ws.Range("K126").Value = 6600
$ws.Range("L126").Value = 10854.546
$ws.Range("M126").Value = -4130
$ws.Range("N126").Value = -15794.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1583.25
$ws.Range("I93").Value = 1424.875
$ws.Range("K93").Value = 1424.875
$ws.Range("M93").Value = -176.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4990.6665
$ws.Range("I100").Value = 5488.3335
$ws.Range("K100").Value = 5488.3335
$ws.Range("M100").Value = -4947.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3401.5247
$ws.Range("I136").Value = 3299.8704
$ws.Range("K136").Value = 9899.611199999999
$ws.Range("M136").Value = -7349.611199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1818.4318
$ws.Range("I136").Value = 1678.2222
$ws.Range("K136").Value = 5034.6666
$ws.Range("M136").Value = -2484.6666
